$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order the "Usage" columns (AF:AJ) on row 2 so that each data value
# follows its semantic header (pkmUsage, tkm-N1Usage, tkm-N2Usage,
# tkm-SZMUsage, tkm-N3Usage), matching the new header order introduced by
# the shared-string table re-sort in the target workbook.
$ws.Range("AF1").Value = "pkmUsage"
$ws.Range("AG1").Value = "tkm-N1Usage"
$ws.Range("AH1").Value = "tkm-N2Usage"
$ws.Range("AI1").Value = "tkm-SZMUsage"
$ws.Range("AJ1").Value = "tkm-N3Usage"

$ws.Range("AF2").Value = 858
$ws.Range("AG2").Value = 7.5
$ws.Range("AH2").Value = 24.2
$ws.Range("AI2").Value = 414.5
$ws.Range("AJ2").Value = 130.3

# --- Small floating point refinements (re-derived recalculation values)
$ws.Range("B2").Value = 42.5073980631569
$ws.Range("F2").Value = 612.0528792959249
$ws.Range("N2").Value = 514.1244186085769
$ws.Range("O2").Value = 514.1244186085769
$ws.Range("T2").Value = 569.5603395168249
